$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trabajo")

# Copy row 35 (values + formatting) down into the new row 36, then
# overwrite the handful of cells that actually change for the new entry.
$ws.Range("A35:K35").Copy($ws.Range("A36:K36"))
$ws.Range("A36").EntireRow.RowHeight = 120

$ws.Cells.Item(36, 2).Value = 35
$ws.Cells.Item(36, 5).Value = "https://www.mitradel.gob.pa/presentan-proyecto-de-ley-que-amplia-el-alcance-de-los-vales-alimenticios/"
$ws.Cells.Item(36, 6).Value = "El Órgano Ejecutivo, a través de la ministra de Trabajo y Desarrollo Laboral, Doris Zapata Acevedo, presentó este miércoles 5 de agosto, de forma virtual, ante el Pleno de la Asamblea Nacional de Diputados, la propuesta de modificación de la Ley 59 del 7 de agosto de 2003, sobre el Programa de Alimentación de Trabajadores, que son los incentivos de productividad que ofrecen los empleadores a sus trabajadores, sin que esto se considere parte del salario."
$ws.Cells.Item(36, 8).Value = 44048
$ws.Cells.Item(36, 9).Value = 44048

# Extend the Excel table / autofilter to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K36"))

# Extend the category data-validation from C2:C35 to C2:C36
$ws.Range("C2:C35").Validation.Delete()
$val = $ws.Range("C2:C36").Validation
$val.Add(0, 1, 1, "")
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true
$val.InputTitle = "Categoria"
$val.InputMessage = "Selecciona una categoría de la lista"
$val.ErrorTitle = "Entrada no válida"
$val.ErrorMessage = "Selecciona una categoría de la lista"

# New hyperlinks for the added row
$ws.Hyperlinks.Add($ws.Range("G36"), "https://www.mitradel.gob.pa/")
$ws.Hyperlinks.Add($ws.Range("E36"), "https://www.mitradel.gob.pa/presentan-proyecto-de-ley-que-amplia-el-alcance-de-los-vales-alimenticios/")

# Adding hyperlinks re-styles the touched cells as generic "Hyperlink"
# style; restore the table's normal look (matching column E/G elsewhere).
$ws.Range("E35").Copy()
$ws.Range("E36").PasteSpecial(-4122)
$ws.Range("G35").Copy()
$ws.Range("G36").PasteSpecial(-4122)

# Move the view down to the newly-added row, like the saved workbook
$ws.Range("E36").Select()
